# Auto-generated edit script applying the cryptos.xlsx cell-value update
# described by the commit diff (106 cell changes across rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '91.194.20'
Set-TextValue $ws.Cells.Item(2, 5) '  +1.71%  '
Set-TextValue $ws.Cells.Item(3, 4) '3.144.79'
Set-TextValue $ws.Cells.Item(3, 5) '  +2.32%  '
Set-TextValue $ws.Cells.Item(4, 5) '  +0.66%  '
Set-TextValue $ws.Cells.Item(5, 4) '237.96'
Set-TextValue $ws.Cells.Item(5, 5) '  +0.81%  '
Set-TextValue $ws.Cells.Item(6, 4) '617.24'
Set-TextValue $ws.Cells.Item(6, 5) '  -0.08%  '
Set-TextValue $ws.Cells.Item(7, 4) '1.11'
Set-TextValue $ws.Cells.Item(7, 5) '  +5.86%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.372'
Set-TextValue $ws.Cells.Item(8, 5) '  +3.07%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.999'
Set-TextValue $ws.Cells.Item(9, 5) '  -0.10%  '
Set-TextValue $ws.Cells.Item(10, 4) '3.139.17'
Set-TextValue $ws.Cells.Item(10, 5) '  +2.17%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.736'
Set-TextValue $ws.Cells.Item(11, 5) '  +3.26%  '
Set-TextValue $ws.Cells.Item(12, 5) '  +2.50%  '
Set-TextValue $ws.Cells.Item(13, 5) '  -1.98%  '
Set-TextValue $ws.Cells.Item(14, 4) '35.04'
Set-TextValue $ws.Cells.Item(14, 5) '  -0.41%  '
Set-TextValue $ws.Cells.Item(15, 4) '5.52'
Set-TextValue $ws.Cells.Item(15, 5) '  +2.95%  '
Set-TextValue $ws.Cells.Item(16, 4) '91.405.23'
Set-TextValue $ws.Cells.Item(16, 5) '  +2.00%  '
Set-TextValue $ws.Cells.Item(17, 4) '3.730.50'
Set-TextValue $ws.Cells.Item(17, 5) '  +1.65%  '
Set-TextValue $ws.Cells.Item(18, 4) '3.130.56'
Set-TextValue $ws.Cells.Item(18, 5) '  +0.72%  '
Set-TextValue $ws.Cells.Item(19, 5) '  -1.93%  '
Set-TextValue $ws.Cells.Item(20, 4) '15.06'
Set-TextValue $ws.Cells.Item(20, 5) '  +9.47%  '
Set-TextValue $ws.Cells.Item(21, 4) '5.82'
Set-TextValue $ws.Cells.Item(21, 5) '  +7.89%  '
Set-TextValue $ws.Cells.Item(22, 5) '  -4.68%  '
Set-TextValue $ws.Cells.Item(23, 4) '440.66'
Set-TextValue $ws.Cells.Item(23, 5) '  +1.83%  '
Set-TextValue $ws.Cells.Item(24, 4) '9.09'
Set-TextValue $ws.Cells.Item(24, 5) '  +3.88%  '
Set-TextValue $ws.Cells.Item(25, 4) '5.71'
Set-TextValue $ws.Cells.Item(25, 5) '  +2.39%  '
Set-TextValue $ws.Cells.Item(26, 2) 'Aptos'
Set-TextValue $ws.Cells.Item(26, 3) 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Cells.Item(26, 4) '11.86'
Set-TextValue $ws.Cells.Item(26, 5) '  +0.96%  '
Set-TextValue $ws.Cells.Item(27, 2) 'Litecoin'
Set-TextValue $ws.Cells.Item(27, 3) 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Cells.Item(27, 4) '82.14'
Set-TextValue $ws.Cells.Item(27, 5) '  -5.20%  '
Set-TextValue $ws.Cells.Item(29, 5) '  +0.19%  '
Set-TextValue $ws.Cells.Item(30, 4) '0.231'
Set-TextValue $ws.Cells.Item(30, 5) '  +18.31%  '
Set-TextValue $ws.Cells.Item(31, 5) '  +9.49%  '
Set-TextValue $ws.Cells.Item(32, 5) '  +38.13%  '
Set-TextValue $ws.Cells.Item(33, 4) '9.29'
Set-TextValue $ws.Cells.Item(33, 5) '  +2.69%  '
Set-TextValue $ws.Cells.Item(34, 2) 'Binance-PegBSC-USD'
Set-TextValue $ws.Cells.Item(34, 3) 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Cells.Item(34, 4) '1.00'
Set-TextValue $ws.Cells.Item(34, 5) '  +0.06%  '
Set-TextValue $ws.Cells.Item(35, 2) 'Kaspa'
Set-TextValue $ws.Cells.Item(35, 3) 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Cells.Item(35, 4) '0.170'
Set-TextValue $ws.Cells.Item(35, 5) '  +11.85%  '
Set-TextValue $ws.Cells.Item(36, 2) 'RenderToken'
Set-TextValue $ws.Cells.Item(36, 3) 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Cells.Item(36, 4) '7.58'
Set-TextValue $ws.Cells.Item(36, 5) '  +6.75%  '
Set-TextValue $ws.Cells.Item(37, 2) 'EthereumClassic'
Set-TextValue $ws.Cells.Item(37, 3) 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(37, 4) '26.20'
Set-TextValue $ws.Cells.Item(37, 5) '  +2.58%  '
Set-TextValue $ws.Cells.Item(38, 4) '504.24'
Set-TextValue $ws.Cells.Item(38, 5) '  +1.90%  '
Set-TextValue $ws.Cells.Item(39, 4) '1.35'
Set-TextValue $ws.Cells.Item(39, 5) '  +7.17%  '
Set-TextValue $ws.Cells.Item(40, 5) '  +2.51%  '
Set-TextValue $ws.Cells.Item(41, 5) '  +12.78%  '
Set-TextValue $ws.Cells.Item(42, 4) '3.79'
Set-TextValue $ws.Cells.Item(42, 5) '  +5.58%  '
Set-TextValue $ws.Cells.Item(43, 5) '  -7.48%  '
Set-TextValue $ws.Cells.Item(44, 4) '22.16'
Set-TextValue $ws.Cells.Item(44, 5) '  +0.27%  '
Set-TextValue $ws.Cells.Item(45, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(46, 4) '159.73'
Set-TextValue $ws.Cells.Item(46, 5) '  +5.17%  '
Set-TextValue $ws.Cells.Item(47, 2) 'Stacks'
Set-TextValue $ws.Cells.Item(47, 3) 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(47, 4) '1.92'
Set-TextValue $ws.Cells.Item(47, 5) '  +4.02%  '
Set-TextValue $ws.Cells.Item(48, 2) 'ARBITRUM'
Set-TextValue $ws.Cells.Item(48, 3) 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Cells.Item(48, 4) '0.704'
Set-TextValue $ws.Cells.Item(48, 5) '  +4.56%  '
Set-TextValue $ws.Cells.Item(49, 4) '1.35'
Set-TextValue $ws.Cells.Item(49, 5) '  +4.45%  '
Set-TextValue $ws.Cells.Item(50, 2) 'OKB'
Set-TextValue $ws.Cells.Item(50, 3) 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Cells.Item(50, 4) '44.00'
Set-TextValue $ws.Cells.Item(50, 5) '  -0.88%  '
Set-TextValue $ws.Cells.Item(51, 2) 'Filecoin'
Set-TextValue $ws.Cells.Item(51, 3) 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(51, 4) '4.39'
Set-TextValue $ws.Cells.Item(51, 5) '  +1.48%  '
